$d = $word.ActiveDocument

# Locate the two bullet paragraphs we need to touch by their current text,
# rather than a hard-coded paragraph index, so the script is resilient to
# the exact paragraph numbering.
$pScanner = $null
$pMetody = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Klasa Scanner*" -or $t -like "Klasa *Scanner*") {
        $pScanner = $p
    }
    if ($t -like "Klasa, metody, pola i obiekty*") {
        $pMetody = $p
    }
}

# --- Paragraph: "Klasa Scanner " ---
# Merge the existing split runs ("Klasa " / "Scanner" / " ") into a single run
# and drop the spell-check proofErr markers around "Scanner" by doing a plain
# (non-tracked) text replace of the full visible text with itself.
$d.Content.Find.Execute("Klasa Scanner ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Klasa Scanner ", 2) | Out-Null

# Now append a new, distinct run "- stereotypy" right after it. Wrapping the
# insertion in a tracked-changes session and then accepting it forces the
# engine to keep the new text as its own <w:r> instead of silently merging it
# back into the previous run (which otherwise happens for adjacent runs that
# share identical formatting).
$d.TrackRevisions = $true
$r1 = $pScanner.Range
$ip1 = $d.Range($r1.End - 1, $r1.End - 1)
$ip1.InsertAfter("- stereotypy")
$d.TrackRevisions = $false
$d.Revisions.AcceptAll() | Out-Null

# --- Paragraph: "Klasa, metody, pola i obiekty" ---
# Append a new run ", modyfikator dostępu" using the same technique so it
# remains a separate <w:r> from the existing text.
$d.TrackRevisions = $true
$r2 = $pMetody.Range
$ip2 = $d.Range($r2.End - 1, $r2.End - 1)
$ip2.InsertAfter(", modyfikator dostępu")
$d.TrackRevisions = $false
$d.Revisions.AcceptAll() | Out-Null

Write-Host "done"
